$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-31 from 45212 to 45221
$ws.Range("C2:C31").Value = 45221
